$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "duration" column before column D (current: label, filename,
# extension, projectionscheme, hmd -> after: label, filename, duration,
# extension, projectionscheme, hmd)
$ws.Range("C1").EntireColumn.Insert()

$ws.Range("C1").Value = "duration"
$ws.Range("C1").Font.Bold = $true

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = 11
}

# All hmd values become "vive" (the "rift" string is retired entirely)
$ws.Range("F2").Value = "vive"
$ws.Range("F4").Value = "vive"
$ws.Range("F6").Value = "vive"
$ws.Range("F8").Value = "vive"

$ws.Range("F11").Select()
